$wb = $excel.ActiveWorkbook

# The "measurement" sheet (Sheet1) has an extra column (M) that duplicates
# data now carried by what was column N. Remove column M so the old column
# N shifts left and becomes the new column M.
$ws = $wb.Worksheets.Item(1)
$ws.Columns("M").Delete()

# Match the author's resulting selection on the sheet (now pointing at the
# former column N, which has become column M).
$ws.Range("M1").Select()
